$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the value is stored as a string (matches source data).
$textForceCells = @("D4","D5","D6","D8","D9","D11","D12","D13","D16","D17","D19","D21","D22","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.051.18"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").Value = "1.805.09"
$ws.Range("E3").Value = "  +4.25%  "
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "315.69"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +12.32%  "
$ws.Range("D8").Value = "0.3804"
$ws.Range("E8").Value = "  +8.21%  "
$ws.Range("D9").Value = "43.11"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  +4.17%  "
$ws.Range("D11").Value = "1.136"
$ws.Range("E11").Value = "  +8.04%  "
$ws.Range("D12").Value = "0.9981"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "21.16"
$ws.Range("E13").Value = "  +5.52%  "
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").Value = "1.801.54"
$ws.Range("E15").Value = "  +4.02%  "
$ws.Range("D16").Value = "7.162"
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").Value = "91.98"
$ws.Range("E17").Value = "  +5.29%  "
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "0.06497"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "17.18"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").Value = "5.991"
$ws.Range("E22").Value = "  +5.01%  "
$ws.Range("D23").Value = "28.070.01"
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("D24").Value = "11.24"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "2.090"
$ws.Range("D26").Value = "20.62"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").Value = "155.89"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "2.389"
$ws.Range("D29").Value = "2.009.20"
$ws.Range("E29").Value = "  +4.19%  "
$ws.Range("D30").Value = "123.00"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "1.153"
$ws.Range("E31").Value = "  +9.35%  "
$ws.Range("D32").Value = "0.1037"
$ws.Range("E32").Value = "  +10.87%  "
$ws.Range("D33").Value = "5.751"
$ws.Range("E33").Value = "  +6.38%  "
$ws.Range("D34").Value = "3.600"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").Value = "0.02305"
$ws.Range("E35").Value = "  +5.16%  "
$ws.Range("D36").Value = "0.2127"
$ws.Range("E36").Value = "  +6.31%  "
$ws.Range("D37").Value = "8.666"
$ws.Range("E37").Value = "  +15.35%  "
$ws.Range("D38").Value = "11.53"
$ws.Range("E38").Value = "  +4.62%  "
$ws.Range("D39").Value = "5.023"
$ws.Range("E39").Value = "  +4.89%  "
$ws.Range("D40").Value = "0.06046"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "0.6300"
$ws.Range("E41").Value = "  +4.61%  "
$ws.Range("D42").Value = "0.9983"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.152"
$ws.Range("E43").Value = "  +4.79%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.396"
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").Value = "13.46"
$ws.Range("E45").Value = "  +4.47%  "
$ws.Range("D46").Value = "0.5926"
$ws.Range("E46").Value = "  +4.57%  "
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").Value = "121.92"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "0.06790"
$ws.Range("E51").Value = "  +2.17%  "
